$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (C) column for all existing data rows (2-391)
#    from 45188 to 45189 (date shifted by one day).
$ws.Range("C2:C391").Value = 45189

# 2) Row 391 gains an explicit row height (ht="15" customHeight="1"),
#    matching the rows above it.
$ws.Rows.Item(391).RowHeight = 15

# 3) Append a new record as row 392.
$ws.Range("A392").Value = "A 44138-2023"
$ws.Range("B392").Value = 45188
$ws.Range("C392").Value = 45189
$ws.Range("D392").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E392").Value = "MALÅ"
$ws.Range("F392").Value = "Sveaskog"
$ws.Range("G392").Value = 13.6
$ws.Range("H392").Value = 0
$ws.Range("I392").Value = 0
$ws.Range("J392").Value = 0
$ws.Range("K392").Value = 0
$ws.Range("L392").Value = 0
$ws.Range("M392").Value = 0
$ws.Range("N392").Value = 0
$ws.Range("O392").Value = 0
$ws.Range("P392").Value = 0
$ws.Range("Q392").Value = 0
$ws.Range("R392").Value = ""

# Match the date formatting / wrap-text styling used by the rest of the table
$ws.Range("B392:C392").NumberFormat = "YYYY-MM-DD"
$ws.Range("R392").WrapText = $true
